$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.90532904932644
$ws.Range("C2").Value = 10.58136980179555
$ws.Range("D2").Value = 4.729078971557712
$ws.Range("E2").Value = 11.2443039230366
$ws.Range("F2").Value = 23.43303130976401
$ws.Range("I2").Value = 20.96508356542171
$ws.Range("L2").Value = 9.841459993918527
$ws.Range("O2").Value = 20.96543010837522

# Row 3
$ws.Range("B3").Value = 16.20083200353795
$ws.Range("C3").Value = 10.26389082591402
$ws.Range("D3").Value = 4.684323851077968
$ws.Range("E3").Value = 11.29889216155078
$ws.Range("F3").Value = 23.48647448217628
$ws.Range("I3").Value = 21.12187738033639
$ws.Range("L3").Value = 9.807471450699259
$ws.Range("O3").Value = 21.06552567469965

# Row 4
$ws.Range("B4").Value = 15.75268506531014
$ws.Range("C4").Value = 10.0630030416044
$ws.Range("D4").Value = 4.656487039072693
$ws.Range("E4").Value = 11.33453456465055
$ws.Range("F4").Value = 23.52808438401497
$ws.Range("I4").Value = 21.22394490102226
$ws.Range("L4").Value = 9.788220897778718
$ws.Range("O4").Value = 21.13363487981603

# Row 5
$ws.Range("B5").Value = 15.56637990378482
$ws.Range("C5").Value = 9.979733876691101
$ws.Range("D5").Value = 4.645060004190032
$ws.Range("E5").Value = 11.34959383497547
$ws.Range("F5").Value = 23.54724161465532
$ws.Range("I5").Value = 21.26699470959757
$ws.Range("L5").Value = 9.780788611670516
$ws.Range("O5").Value = 21.16305416718479

# Row 6
$ws.Range("B6").Value = 15.5352292393528
$ws.Range("C6").Value = 9.965825018623475
$ws.Range("D6").Value = 4.643157711752456
$ws.Range("E6").Value = 11.3521267207185
$ws.Range("F6").Value = 23.5505552432968
$ws.Range("I6").Value = 21.27423104176428
$ws.Range("L6").Value = 9.779579550104538
$ws.Range("O6").Value = 21.16803949636331

# Row 7
$ws.Range("B7").Value = 15.75018706004028
$ws.Range("C7").Value = 10.06188560973205
$ws.Range("D7").Value = 4.656333258840439
$ws.Range("E7").Value = 11.33473549376842
$ws.Range("F7").Value = 23.52833384922659
$ws.Range("I7").Value = 21.22451958969522
$ws.Range("L7").Value = 9.788118986372497
$ws.Range("O7").Value = 21.13402491090552

# Row 8
$ws.Range("B8").Value = 16.66578540602316
$ws.Range("C8").Value = 10.47319138236871
$ws.Range("D8").Value = 4.713725355615991
$ws.Range("E8").Value = 11.26268511315002
$ws.Range("F8").Value = 23.44962683290571
$ws.Range("I8").Value = 21.01794268179244
$ws.Range("L8").Value = 9.829408220527997
$ws.Range("O8").Value = 20.99855806804023

# Row 9
$ws.Range("B9").Value = 18.32871560163884
$ws.Range("C9").Value = 11.22888965123296
$ws.Range("D9").Value = 4.823159012711213
$ws.Range("E9").Value = 11.13823946873631
$ws.Range("F9").Value = 23.3654988648326
$ws.Range("I9").Value = 20.65887061883323
$ws.Range("L9").Value = 9.922960637898093
$ws.Range("O9").Value = 20.78602393239085

# Row 10
$ws.Range("B10").Value = 19.45960672573077
$ws.Range("C10").Value = 11.74877091343444
$ws.Range("D10").Value = 4.901295433160184
$ws.Range("E10").Value = 11.05705609475029
$ws.Range("F10").Value = 23.34699124613474
$ws.Range("I10").Value = 20.42317635572227
$ws.Range("L10").Value = 9.998994957978043
$ws.Range("O10").Value = 20.66272354443375

# Row 11
$ws.Range("B11").Value = 19.9526972741925
$ws.Range("C11").Value = 11.97683905555715
$ws.Range("D11").Value = 4.936274066615988
$ws.Range("E11").Value = 11.02234437297453
$ws.Range("F11").Value = 23.34804994690561
$ws.Range("I11").Value = 20.32207362663525
$ws.Range("L11").Value = 10.03508375501266
$ws.Range("O11").Value = 20.61385969187156

# Row 12
$ws.Range("B12").Value = 20.13623565170918
$ws.Range("C12").Value = 12.06193542354952
$ws.Range("D12").Value = 4.94943192342184
$ws.Range("E12").Value = 11.00951880747862
$ws.Range("F12").Value = 23.34981750186194
$ws.Range("I12").Value = 20.28466989770184
$ws.Range("L12").Value = 10.04895753269622
$ws.Range("O12").Value = 20.5964028822474

# Row 13
$ws.Range("B13").Value = 20.09685061905292
$ws.Range("C13").Value = 12.04366562360644
$ws.Range("D13").Value = 4.946602148412051
$ws.Range("E13").Value = 11.01226683542257
$ws.Range("F13").Value = 23.34937600292955
$ws.Range("I13").Value = 20.29268620941857
$ws.Range("L13").Value = 10.04596045842613
$ws.Range("O13").Value = 20.60011583799426

# Row 14
$ws.Range("B14").Value = 19.96786149169695
$ws.Range("C14").Value = 11.98386569843754
$ws.Range("D14").Value = 4.937358367026524
$ws.Range("E14").Value = 11.02128281290873
$ws.Range("F14").Value = 23.34816795939212
$ws.Range("I14").Value = 20.31897871569379
$ws.Range("L14").Value = 10.03622104280502
$ws.Range("O14").Value = 20.61240248529739

# Row 15
$ws.Range("B15").Value = 19.88843417511351
$ws.Range("C15").Value = 11.94706981723114
$ws.Range("D15").Value = 4.931684663476753
$ws.Range("E15").Value = 11.0268469020629
$ws.Range("F15").Value = 23.34760605911133
$ws.Range("I15").Value = 20.33519852529238
$ws.Range("L15").Value = 10.03028217919829
$ws.Range("O15").Value = 20.62006497075793

# Row 16
$ws.Range("B16").Value = 19.4269422568432
$ws.Range("C16").Value = 11.73369158361997
$ws.Range("D16").Value = 4.898997557000831
$ws.Range("E16").Value = 11.05936923201921
$ws.Range("F16").Value = 23.3471131333948
$ws.Range("I16").Value = 20.42990690922032
$ws.Range("L16").Value = 9.996665983489869
$ws.Range("O16").Value = 20.66606298994536

# Row 17
$ws.Range("B17").Value = 19.13827556553996
$ws.Range("C17").Value = 11.60059029852801
$ws.Range("D17").Value = 4.878795505108132
$ws.Range("E17").Value = 11.07988890312775
$ws.Range("F17").Value = 23.34924142666573
$ws.Range("I17").Value = 20.48957536882275
$ws.Range("L17").Value = 9.976422217531233
$ws.Range("O17").Value = 20.69613767597885

# Row 18
$ws.Range("B18").Value = 18.97023740928912
$ws.Range("C18").Value = 11.52324357080238
$ws.Range("D18").Value = 4.867122916894973
$ws.Range("E18").Value = 11.09190014139607
$ws.Range("F18").Value = 23.35135754760665
$ws.Range("I18").Value = 20.52447058914454
$ws.Range("L18").Value = 9.964920264680771
$ws.Range("O18").Value = 20.71411554787137

# Row 19
$ws.Range("B19").Value = 18.9130020046141
$ws.Range("C19").Value = 11.49692136033857
$ws.Range("D19").Value = 4.863161890517233
$ws.Range("E19").Value = 11.09600281610134
$ws.Range("F19").Value = 23.35222708863806
$ws.Range("I19").Value = 20.53638428864278
$ws.Range("L19").Value = 9.961050481943699
$ws.Range("O19").Value = 20.72031904589739

# Row 20
$ws.Range("B20").Value = 19.16921298280146
$ws.Range("C20").Value = 11.61484138604161
$ws.Range("D20").Value = 4.880951570036138
$ws.Range("E20").Value = 11.07768293289842
$ws.Range("F20").Value = 23.34892251712662
$ws.Range("I20").Value = 20.48316398159914
$ws.Range("L20").Value = 9.978562590008091
$ws.Range("O20").Value = 20.69286576994518

# Row 21
$ws.Range("B21").Value = 20.00583596207962
$ws.Range("C21").Value = 12.00146521764904
$ws.Range("D21").Value = 4.940075922607199
$ws.Range("E21").Value = 11.01862594443466
$ws.Range("F21").Value = 23.34848567885335
$ws.Range("I21").Value = 20.31123202197102
$ws.Range("L21").Value = 10.03907617288773
$ws.Range("O21").Value = 20.60876512686886

# Row 22
$ws.Range("B22").Value = 20.53401616561887
$ws.Range("C22").Value = 12.24673314597908
$ws.Range("D22").Value = 4.97820252418012
$ws.Range("E22").Value = 10.98188820513242
$ws.Range("F22").Value = 23.35616663698186
$ws.Range("I22").Value = 20.20400482268783
$ws.Range("L22").Value = 10.07983219591624
$ws.Range("O22").Value = 20.55990586128441

# Row 23
$ws.Range("B23").Value = 20.25385145608345
$ws.Range("C23").Value = 12.11652427415692
$ws.Range("D23").Value = 4.95790279797651
$ws.Range("E23").Value = 11.0013257109947
$ws.Range("F23").Value = 23.35133741793356
$ws.Range("I23").Value = 20.2607628748061
$ws.Range("L23").Value = 10.05797223620708
$ws.Range("O23").Value = 20.58542183686675

# Row 24
$ws.Range("B24").Value = 19.15523264591532
$ws.Range("C24").Value = 11.60840103487065
$ws.Range("D24").Value = 4.879976993739231
$ws.Range("E24").Value = 11.078679584795
$ws.Range("F24").Value = 23.34906391621472
$ws.Range("I24").Value = 20.48606072886169
$ws.Range("L24").Value = 9.977594501959382
$ws.Range("O24").Value = 20.6943428575898

# Row 25
$ws.Range("B25").Value = 17.89417387060803
$ws.Range("C25").Value = 11.03037386063999
$ws.Range("D25").Value = 4.793925753461214
$ws.Range("E25").Value = 11.17010434038284
$ws.Range("F25").Value = 23.38068381269744
$ws.Range("I25").Value = 20.75107646186897
$ws.Range("L25").Value = 9.896341693515618
$ws.Range("O25").Value = 20.83778455662568
